$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer`" to be shown on Friday.`n"
$ws.Range("D2").Value = "Oppenheimer_was_selected, "
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded with no agreement on a movie to be shown on Friday.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision has been recorded as a no-decision regarding the movie to show on Friday.`n"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been successfully acquired for showing on Friday.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision was made not to acquire any movie for Friday, and the conversation ended without a consensus on a specific film.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The function for no decision has been executed successfully.`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been made, and no specific movie will be shown on Friday.`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday is that there was no definite selection made.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: No decision has been made regarding the movie for Friday.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
$ws.Range("D11").Value = "both_movies, "
$ws.Range("C12").Value = "MSG: None`n`nMSG: The committee ended the conversation without a decision regarding which movie to show on Friday.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made, as the committee could not reach a consensus.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision has been successfully recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The conversation has concluded without reaching a decision about which movie to show on Friday.`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement was reached regarding the movie for Friday.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision has been recorded with no agreement on which movie to show on Friday.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: I have recorded the decision: there was no decision made regarding the movie to be shown on Friday.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: I have recorded the decision as no decision made regarding the movie to show on Friday.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: None`n`nMSG: no_decision`n"
$ws.Range("D21").Value = "no_decision, , no_decision, "
$ws.Range("C22").Value = "MSG: None`n`nMSG: No decision was made regarding the movie to show on Friday.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision is recorded as no decision was made regarding which movie to show on Friday.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision-making process for Friday's movie has concluded without a clear choice being made.`n"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision process concluded with no movie selected for Friday.`n"
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie`" for the movie to be shown on Friday.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was made regarding which movie to show on Friday.`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision has been recorded. No movie was selected for Friday's showing.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for both movies.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.`n"
$ws.Range("D33").Value = "Barbie_was_selected, "
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision was recorded, indicating that there is no agreement on a movie to show on Friday.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" will be the movie shown on Friday.`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision process concluded without a clear choice of a movie for Friday.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The conversation ended without a clear decision on which movie to show on Friday. Therefore, I will call the no_decision function.`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision being made regarding the movie for Friday.`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: No decision was made regarding the movie to be shown on Friday.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision to show a movie on Friday was not made, resulting in no selection.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision-making process ended without a clear choice for a movie to show on Friday.`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision to acquire rights for both movies has been recorded.`n"
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision about which movie to play on Friday was not finalized.`n"
$ws.Range("C47").Value = "MSG: None`n`nMSG: No decision was made regarding the movie to show on Friday.`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer`" to be shown on Friday.`n"
$ws.Range("C49").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been successfully selected for Friday's showing.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for the movie `"Barbie.`"`n"
$ws.Range("C51").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision has been recorded as no specific movie was chosen.`n"
$ws.Range("D52").Value = "no_decision, "
$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding Friday's movie was made.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision has been recorded with no definitive choice made regarding the movie for Friday.`n"
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision about Friday's movie was made.`n"
$ws.Range("C56").Value = "MSG: None`n`nMSG: The function for no decision has been successfully called, indicating that a consensus on the movie to be shown on Friday was not reached.`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has resulted in no conclusive agreement.`n"
$ws.Range("C59").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired successfully.`n"
$ws.Range("C60").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made.`n"
$ws.Range("C61").Value = "MSG: None`n`nMSG: The decision process has concluded without a specific agreement on which movie to show on Friday.`n"
$ws.Range("C62").Value = "MSG: None`n`nMSG: The rights to both movies, `"Oppenheimer`" and `"Barbie,`" have been acquired.`n"
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision resulted in no agreement on which movie to show on Friday.`n"
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both `"Oppenheimer`" and `"Barbie`" has been successfully recorded.`n"
$ws.Range("D64").Value = "both_movies, "
$ws.Range("C65").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday has concluded without a selection being made. If you need assistance with anything else, feel free to ask!`n"
$ws.Range("C66").Value = "MSG: None`n`nMSG: The decision has been made to select `"Barbie`" for the movie to be shown on Friday.`n"
$ws.Range("C67").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has not been finalized.`n"
$ws.Range("D67").Value = "no_decision, "

Write-Host "applied updates"
